$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Sheet1 ("Excel Sheet Model" gains an Age column)
# ---------------------------------------------------------------------
$ws1.Range("D2:D3").NumberFormat = "0"
$ws1.Cells.Item(2, 4).Value = 34
$ws1.Cells.Item(3, 4).Value = 28

$ws1.Cells.Item(1, 4).Value = "Age"
$ws1.Cells.Item(1, 4).Font.Bold = $true
$ws1.Cells.Item(1, 4).NumberFormat = "0"

# ---------------------------------------------------------------------
# Sheet2 (new gym-visit rows sent, one per JDBC row, via the Camel route)
# ---------------------------------------------------------------------
$ws2.Cells.Item(1, 1).Value = "MemberID"
$ws2.Cells.Item(1, 1).Font.Bold = $true
$ws2.Cells.Item(1, 1).NumberFormat = "0"
$ws2.Cells.Item(1, 2).Value = "VisitingDate"
$ws2.Cells.Item(1, 2).Font.Bold = $true
$ws2.Cells.Item(1, 2).NumberFormat = "m/dd/yy;@"
$ws2.Cells.Item(1, 3).Value = "ExerciseZone"
$ws2.Cells.Item(1, 3).Font.Bold = $true

$ws2.Cells.Item(3, 2).NumberFormat = "@"
$ws2.Cells.Item(3, 2).Value = "6/19/2017"
$ws2.Cells.Item(3, 2).NumberFormat = "dd/mm/yyyy;@"

$ws2.Cells.Item(2, 1).NumberFormat = "0"
$ws2.Cells.Item(2, 1).Value = 6
$ws2.Cells.Item(2, 2).NumberFormat = "@"
$ws2.Cells.Item(2, 2).Value = "5/14/2017"
$ws2.Cells.Item(2, 2).NumberFormat = "dd/mm/yyyy;@"
$ws2.Cells.Item(2, 2).WrapText = $false
$ws2.Cells.Item(2, 3).Value = "Track"

$ws2.Cells.Item(3, 1).NumberFormat = "0"
$ws2.Cells.Item(3, 1).Value = 7
$ws2.Cells.Item(3, 3).Value = "HeavyLifting"

$ws2.Columns.Item(1).AutoFit()
$ws2.Columns.Item(2).AutoFit()
$ws2.Columns.Item(3).AutoFit()

$ws2.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping
# ---------------------------------------------------------------------
$ws1.Range("F10").Select()
$ws2.Activate()
$ws2.Range("B4").Select()
